$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1734.6
$ws.Range("I40").Value = 1768.625
$ws.Range("J40").Value = 1705.9474
$ws.Range("K40").Value = 1768.625
$ws.Range("L40").Value = 1705.9474
$ws.Range("M40").Value = -1593.625
$ws.Range("N40").Value = -2055.9474

$ws.Range("H107").Value = 581.7143
$ws.Range("I107").Value = 624.64514
$ws.Range("J107").Value = 249
$ws.Range("K107").Value = 624.64514
$ws.Range("L107").Value = 249
$ws.Range("M107").Value = 1295.35486
$ws.Range("N107").Value = -4089

$ws.Range("H113").Value = 9632.299999999999
$ws.Range("I113").Value = 19855.363
$ws.Range("J113").Value = 3713.6843
$ws.Range("K113").Value = 19855.363
$ws.Range("L113").Value = 3713.6843
$ws.Range("M113").Value = -16601.363
$ws.Range("N113").Value = -10221.6843

$ws.Range("H132").Value = 29457632
$ws.Range("I132").Value = 32804892
$ws.Range("J132").Value = 1747.2
$ws.Range("K132").Value = 98414676
$ws.Range("L132").Value = 5241.6
$ws.Range("M132").Value = -98412146
$ws.Range("N132").Value = -10301.6

$ws.Range("H137").Value = 170972.14
$ws.Range("I137").Value = 196002.05
$ws.Range("J137").Value = 1325
$ws.Range("K137").Value = 588006.1499999999
$ws.Range("L137").Value = 3975
$ws.Range("M137").Value = -585456.1499999999
$ws.Range("N137").Value = -9075

$ws.Range("H138").Value = 1954.5747
$ws.Range("I138").Value = 1302.875
$ws.Range("J138").Value = 2509.2126
$ws.Range("K138").Value = 3908.625
$ws.Range("L138").Value = 7527.6378
$ws.Range("M138").Value = 1231.375
$ws.Range("N138").Value = -17807.6378

$ws.Range("H141").Value = 3194.4084
$ws.Range("I141").Value = 2570.07
$ws.Range("J141").Value = 5736.357
$ws.Range("K141").Value = 7710.210000000001
$ws.Range("L141").Value = 17209.071
$ws.Range("M141").Value = -2530.210000000001
$ws.Range("N141").Value = -27569.071

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5541.5796
$ws.Range("I32").Value = 4397.8335
$ws.Range("J32").Value = 14462.8
$ws.Range("K32").Value = 4397.8335
$ws.Range("L32").Value = 14462.8
$ws.Range("M32").Value = -4110.8335
$ws.Range("N32").Value = -15036.8

$ws.Range("H61").Value = 2911.7874
$ws.Range("I61").Value = 2880.5898
$ws.Range("J61").Value = 3063.875
$ws.Range("K61").Value = 2880.5898
$ws.Range("L61").Value = 3063.875
$ws.Range("M61").Value = -2668.5898
$ws.Range("N61").Value = -3487.875

$ws.Range("H74").Value = 683.9838999999999
$ws.Range("I74").Value = 500.07144
$ws.Range("J74").Value = 1070.2
$ws.Range("K74").Value = 500.07144
$ws.Range("L74").Value = 1070.2
$ws.Range("M74").Value = 373.92856
$ws.Range("N74").Value = -2818.2

$ws.Range("H77").Value = 683.9838999999999
$ws.Range("I77").Value = 500.07144
$ws.Range("J77").Value = 1070.2
$ws.Range("K77").Value = 2500.3572
$ws.Range("L77").Value = 5351
$ws.Range("M77").Value = 1867.6428
$ws.Range("N77").Value = -14087

$ws.Range("H122").Value = 22728388
$ws.Range("I122").Value = 30303964
$ws.Range("K122").Value = 90911892
$ws.Range("M122").Value = -90909442

$ws.Range("H132").Value = 3572970.2
$ws.Range("I132").Value = 4630808
$ws.Range("J132").Value = 2767.5
$ws.Range("K132").Value = 13892424
$ws.Range("L132").Value = 8302.5
$ws.Range("M132").Value = -13889894
$ws.Range("N132").Value = -13362.5

$ws.Range("H136").Value = 2911.7874
$ws.Range("I136").Value = 2880.5898
$ws.Range("J136").Value = 3063.875
$ws.Range("K136").Value = 8641.769400000001
$ws.Range("L136").Value = 9191.625
$ws.Range("M136").Value = -6091.769400000001
$ws.Range("N136").Value = -14291.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7949626.5
$ws.Range("I134").Value = 9023660
$ws.Range("J134").Value = 1780
$ws.Range("K134").Value = 27070980
$ws.Range("L134").Value = 5340
$ws.Range("M134").Value = -27068445
$ws.Range("N134").Value = -10410

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5752.2812
$ws.Range("I31").Value = 980.89795
$ws.Range("J31").Value = 21338.8
$ws.Range("K31").Value = 980.89795
$ws.Range("L31").Value = 21338.8
$ws.Range("M31").Value = -685.89795
$ws.Range("N31").Value = -21928.8

$ws.Range("H34").Value = 5752.2812
$ws.Range("I34").Value = 980.89795
$ws.Range("J34").Value = 21338.8
$ws.Range("K34").Value = 980.89795
$ws.Range("L34").Value = 21338.8
$ws.Range("M34").Value = -778.89795
$ws.Range("N34").Value = -21742.8

$ws.Range("H58").Value = 3348884.8
$ws.Range("I58").Value = 4796067.5
$ws.Range("J58").Value = 9231.538
$ws.Range("K58").Value = 4796067.5
$ws.Range("L58").Value = 9231.538
$ws.Range("M58").Value = -4795864.5
$ws.Range("N58").Value = -9637.538

$ws.Range("H136").Value = 3348884.8
$ws.Range("I136").Value = 4796067.5
$ws.Range("J136").Value = 9231.538
$ws.Range("K136").Value = 14388202.5
$ws.Range("L136").Value = 27694.614
$ws.Range("M136").Value = -14385652.5
$ws.Range("N136").Value = -32794.614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 20042550
$ws.Range("I137").Value = 57933.277
$ws.Range("J137").Value = 71431570
$ws.Range("K137").Value = 173799.831
$ws.Range("L137").Value = 214294710
$ws.Range("M137").Value = -168699.831
$ws.Range("N137").Value = -214304910

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 23810386
$ws.Range("I102").Value = 27778618
$ws.Range("J102").Value = 1001.3333
$ws.Range("K102").Value = 27778618
$ws.Range("L102").Value = 1001.3333
$ws.Range("M102").Value = -27776996
$ws.Range("N102").Value = -4245.3333

$ws.Range("H113").Value = 1516.6666
$ws.Range("I113").Value = 1375
$ws.Range("J113").Value = 1587.5
$ws.Range("K113").Value = 1375
$ws.Range("L113").Value = 1587.5
$ws.Range("M113").Value = 795
$ws.Range("N113").Value = -5927.5

$ws.Range("H132").Value = 17244172
$ws.Range("I132").Value = 25642410
$ws.Range("J132").Value = 5682.421
$ws.Range("K132").Value = 76927230
$ws.Range("L132").Value = 17047.263
$ws.Range("M132").Value = -76924700
$ws.Range("N132").Value = -22107.263

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1408.0714
$ws.Range("I7").Value = 1411.0834
$ws.Range("J7").Value = 1390
$ws.Range("K7").Value = 1411.0834
$ws.Range("L7").Value = 1390
$ws.Range("M7").Value = -1299.0834
$ws.Range("N7").Value = -1614

$ws.Range("H126").Value = 1408.0714
$ws.Range("I126").Value = 1411.0834
$ws.Range("J126").Value = 1390
$ws.Range("K126").Value = 4233.2502
$ws.Range("L126").Value = 4170
$ws.Range("M126").Value = -1763.2502
$ws.Range("N126").Value = -9110

$ws.Range("H136").Value = 3278.049
$ws.Range("I136").Value = 3443.12
$ws.Range("J136").Value = 2527.7273
$ws.Range("K136").Value = 10329.36
$ws.Range("L136").Value = 7583.1819
$ws.Range("M136").Value = -7779.360000000001
$ws.Range("N136").Value = -12683.1819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7934693
$ws.Range("I132").Value = 4000900.8
$ws.Range("J132").Value = 12851934
$ws.Range("K132").Value = 12002702.4
$ws.Range("L132").Value = 38555802
$ws.Range("M132").Value = -12000172.4
$ws.Range("N132").Value = -38560862

$ws.Range("H136").Value = 11904802
$ws.Range("I136").Value = 6867704.5
$ws.Range("J136").Value = 25001258
$ws.Range("K136").Value = 20603113.5
$ws.Range("L136").Value = 75003774
$ws.Range("M136").Value = -20600563.5
$ws.Range("N136").Value = -75008874
